$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "Investor"
$ws.Range("E2").Value = "Kalaari Capital"
$ws.Range("E3").Value = "Accel"

$ws.Range("E4").Select()
